$wb = $excel.ActiveWorkbook
$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux = $wb.Worksheets.Item("Totaux")

# --- Journal sheet: fill in row 10 (C), add new row 11, add blank styled row 12 ---

# Row 10: set the missing "Temps [h]" value (1h05 = 65 minutes)
$wsJournal.Range("C10").Value = 0.045138888888888888

# Row 11: new entry - 07/02/2023, week 2, 25 minutes, Documentation, Redaction du rapport de projet
$wsJournal.Range("A11").Value = 44964
$wsJournal.Range("B11").Value = 2
$wsJournal.Range("C11").Value = 0.017361111111111112
$wsJournal.Range("D11").Value = "Documentation"
$wsJournal.Range("E11").Value = "Rédaction du rapport de projet"

# Resize Tableau1 to include the new row, leaving one blank templated row after it
$wsJournal.ListObjects.Item("Tableau1").Resize($wsJournal.Range("A1:E12"))

# Row 12: blank templated row below, carrying over the same cell formatting as row 11
$wsJournal.Range("A11:D11").Copy()
$wsJournal.Range("A12:D12").PasteSpecial(-4122)

# --- Totaux sheet: add new weekly total row for week of 06/02/2023 ---
$wsTotaux.Range("A4").Value = 44964
$wsTotaux.Range("B4").Formula = "=SUM(Journal!C10:C11)"

$wsTotaux.ListObjects.Item("Tableau2").Resize($wsTotaux.Range("A1:B5"))

# --- View state: active sheet becomes Journal, with selection at E19 ---
$wsTotaux.Range("B9").Select()
$wsJournal.Activate()
$wsJournal.Range("E19").Select()

$wb.Save()
